$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.159.38"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.861.73"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.00"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2844"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.32"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06538"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.35"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07748"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "1.903.90"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.28"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6872"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.075"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "264.71"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "30.150.39"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007727"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "2.103.70"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.218"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.141"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.438"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.96"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.921"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09899"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.327"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.458"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.029"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04715"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6951"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01855"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.761"
$ws.Range("E40").Value = "  +6.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.286"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.84"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.919"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8316"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4118"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.56"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "966.45"
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.083"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.132"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.53"
$ws.Range("E51").Value = "  +1.15%  "
